# Auto-generated edit script applying the Tonberry_Profits price/profit recalculation update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (sheet ALC)
$ws.Range("H2").Value = 388.53845
$ws.Range("I2").Value = 243
$ws.Range("K2").Value = 243
$ws.Range("M2").Value = -130

# Row 40 (sheet ALC)
$ws.Range("H40").Value = 2800
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2800
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3150

# Row 100 (sheet ALC)
$ws.Range("H100").Value = 833.8461
$ws.Range("I100").Value = 820
$ws.Range("K100").Value = 820
$ws.Range("M100").Value = -279

# Row 116 (sheet ALC)
$ws.Range("H116").Value = 8874.375
$ws.Range("I116").Value = 11535.637
$ws.Range("J116").Value = 3019.6
$ws.Range("K116").Value = 11535.637
$ws.Range("L116").Value = 3019.6
$ws.Range("M116").Value = -8093.637000000001
$ws.Range("N116").Value = -9903.6

# Row 138 (sheet ALC)
$ws.Range("H138").Value = 2414.0212
$ws.Range("I138").Value = 2698.9443
$ws.Range("J138").Value = 2237.1724
$ws.Range("K138").Value = 8096.8329
$ws.Range("L138").Value = 6711.5172
$ws.Range("M138").Value = -2956.8329
$ws.Range("N138").Value = -16991.5172

$ws = $wb.Worksheets.Item("ARM")
# Row 16 (sheet ARM)
$ws.Range("H16").Value = 15000
$ws.Range("I16").Value = 10000
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = -9713
$ws.Range("N16").Value = -20574

# Row 102 (sheet ARM)
$ws.Range("H102").Value = 1468.7778
$ws.Range("I102").Value = 1468.7778
$ws.Range("K102").Value = 1468.7778
$ws.Range("M102").Value = 153.2221999999999

# Row 122 (sheet ARM)
$ws.Range("H122").Value = 1270.1428
$ws.Range("I122").Value = 1298.6923
$ws.Range("K122").Value = 3896.0769
$ws.Range("M122").Value = -1446.0769

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (sheet BSM)
$ws.Range("H99").Value = 1518.7273
$ws.Range("I99").Value = 1551
$ws.Range("J99").Value = 1480
$ws.Range("K99").Value = 1551
$ws.Range("L99").Value = 1480
$ws.Range("M99").Value = -53
$ws.Range("N99").Value = -4476

# Row 102 (sheet BSM)
$ws.Range("H102").Value = 1185.3334
$ws.Range("I102").Value = 1185.3334
$ws.Range("K102").Value = 1185.3334
$ws.Range("M102").Value = 2059.6666

# Row 105 (sheet BSM)
$ws.Range("H105").Value = 2691.0557
$ws.Range("I105").Value = 2437.5881
$ws.Range("K105").Value = 2437.5881
$ws.Range("M105").Value = -690.5880999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (sheet CRP)
$ws.Range("H22").Value = 1545.3636
$ws.Range("J22").Value = 1881
$ws.Range("L22").Value = 1881
$ws.Range("N22").Value = -2581

# Row 25 (sheet CRP)
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 31 (sheet CRP)
$ws.Range("H31").Value = 2635.0625
$ws.Range("I31").Value = 2162.8572
$ws.Range("J31").Value = 3002.3333
$ws.Range("K31").Value = 2162.8572
$ws.Range("L31").Value = 3002.3333
$ws.Range("M31").Value = -1867.8572
$ws.Range("N31").Value = -3592.3333

# Row 34 (sheet CRP)
$ws.Range("H34").Value = 2635.0625
$ws.Range("I34").Value = 2162.8572
$ws.Range("J34").Value = 3002.3333
$ws.Range("K34").Value = 2162.8572
$ws.Range("L34").Value = 3002.3333
$ws.Range("M34").Value = -1960.8572
$ws.Range("N34").Value = -3406.3333

# Row 60 (sheet CRP)
$ws.Range("H60").Value = 25862.25
$ws.Range("J60").Value = 25862.25
$ws.Range("L60").Value = 25862.25
$ws.Range("N60").Value = -26884.25

# Row 141 (sheet CRP)
$ws.Range("H141").Value = 64366.5
$ws.Range("J141").Value = 62839.8
$ws.Range("L141").Value = 62839.8
$ws.Range("N141").Value = -73199.8

$ws = $wb.Worksheets.Item("CUL")
# Row 101 (sheet CUL)
$ws.Range("H101").Value = 6166.6665
$ws.Range("J101").Value = 6166.6665
$ws.Range("L101").Value = 18499.9995
$ws.Range("N101").Value = -23367.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 22 (sheet GSM)
$ws.Range("H22").Value = 80009
$ws.Range("J22").Value = 80009
$ws.Range("L22").Value = 80009
$ws.Range("N22").Value = -81067

# Row 27 (sheet GSM)
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 132 (sheet GSM)
$ws.Range("H132").Value = 3209226.8
$ws.Range("I132").Value = 3849572.2
$ws.Range("K132").Value = 11548716.6
$ws.Range("M132").Value = -11546186.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (sheet LTW)
$ws.Range("H22").Value = 1077.8334
$ws.Range("I22").Value = 486.10526
$ws.Range("K22").Value = 486.10526
$ws.Range("M22").Value = -191.10526

# Row 27 (sheet LTW)
$ws.Range("H27").Value = 1077.8334
$ws.Range("I27").Value = 486.10526
$ws.Range("K27").Value = 486.10526
$ws.Range("M27").Value = -379.10526

# Row 38 (sheet LTW)
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 46 (sheet LTW)
$ws.Range("H46").Value = 2381.5625
$ws.Range("I46").Value = 1764
$ws.Range("J46").Value = 2861.889
$ws.Range("K46").Value = 1764
$ws.Range("L46").Value = 2861.889
$ws.Range("M46").Value = -1576
$ws.Range("N46").Value = -3237.889

# Row 80 (sheet LTW)
$ws.Range("H80").Value = 38888
$ws.Range("J80").Value = 38888
$ws.Range("L80").Value = 38888
$ws.Range("N80").Value = -41134

# Row 83 (sheet LTW)
$ws.Range("H83").Value = 38888
$ws.Range("J83").Value = 38888
$ws.Range("L83").Value = 116664
$ws.Range("N83").Value = -127896

# Row 93 (sheet LTW)
$ws.Range("H93").Value = 1240.0454
$ws.Range("I93").Value = 653.94116
$ws.Range("J93").Value = 3232.8
$ws.Range("K93").Value = 653.94116
$ws.Range("L93").Value = 3232.8
$ws.Range("M93").Value = 594.05884
$ws.Range("N93").Value = -5728.8

# Row 100 (sheet LTW)
$ws.Range("H100").Value = 1750
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -3082

# Row 122 (sheet LTW)
$ws.Range("H122").Value = 2723.7
$ws.Range("I122").Value = 2779.75
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 8339.25
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -5889.25
$ws.Range("N122").Value = -12398.5

$ws = $wb.Worksheets.Item("WVR")
# Row 18 (sheet WVR)
$ws.Range("H18").Value = 14125

# Row 21 (sheet WVR)
$ws.Range("H21").Value = 14400
$ws.Range("I21").Value = 13000
$ws.Range("J21").Value = 14750
$ws.Range("K21").Value = 13000
$ws.Range("L21").Value = 14750
$ws.Range("M21").Value = -12765
$ws.Range("N21").Value = -15220

# Row 35 (sheet WVR)
$ws.Range("H35").Value = 14400
$ws.Range("I35").Value = 13000
$ws.Range("J35").Value = 14750
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 14750
$ws.Range("M35").Value = -12710
$ws.Range("N35").Value = -15330

# Row 70 (sheet WVR)
$ws.Range("H70").Value = 46996.5
$ws.Range("J70").Value = 46996.5
$ws.Range("L70").Value = 46996.5
$ws.Range("N70").Value = -47626.5

# Row 73 (sheet WVR)
$ws.Range("H73").Value = 46996.5
$ws.Range("J73").Value = 46996.5
$ws.Range("L73").Value = 46996.5
$ws.Range("N73").Value = -49180.5

# Row 100 (sheet WVR)
$ws.Range("H100").Value = 1298.8334
$ws.Range("I100").Value = 825
$ws.Range("K100").Value = 1650
$ws.Range("M100").Value = -1109
